# Applies the edits described in the commit diff ("update minimum wage
# essay and readme") to the "A Case for AI to Make All of Our Decisions"
# document:
#
#  1. "we're " + "probably" + " not thinking" -> merged into a single run
#     (no visible text change, just a run-split cleanup).
#  2. Drop the proofErr(gramStart/gramEnd) wrapper around "make a decision"
#     in the "Every time we make a decision..." paragraph.
#  3. Rewrite the opening of the "AI Should Make All Our Decisions" section:
#     "Decisions need to be made objectively, otherwise they are hostage to
#     our emotions and biases. Two necessary conditions for an objective
#     decision: the rules and process..." becomes "Our emotions and biases
#     are always going to have an influence on our decisions. For a
#     decision to be objective, the process...".
#  4. Drop the proofErr(gramStart/gramEnd) wrapper around "actually
#     provide"+"d" in the Ethics paragraph.
#  5. Drop the proofErr(gramStart/gramEnd) wrapper around "actually can" in
#     the Ethics paragraph.

$d = $word.ActiveDocument

# Helper: build a <w:p> fragment (all runs sz/szCs 20, matching this essay's
# body-text formatting) out of a list of run-text strings, for use with
# Range.InsertXML - this lets us rewrite a paragraph's runs precisely.
function Build-BodyParagraphXml {
    param([string[]]$RunTexts)

    $xml = "<w:p xmlns:w='http://schemas.openxmlformats.org/wordprocessingml/2006/main'>" +
           "<w:pPr><w:rPr><w:sz w:val='20'/><w:szCs w:val='20'/></w:rPr></w:pPr>"
    foreach ($t in $RunTexts) {
        $escaped = $t -replace '&', '&amp;' -replace '<', '&lt;' -replace '>', '&gt;'
        $xml += "<w:r><w:rPr><w:sz w:val='20'/><w:szCs w:val='20'/></w:rPr>" +
                "<w:t xml:space='preserve'>$escaped</w:t></w:r>"
    }
    $xml += "</w:p>"
    return $xml
}

$rsq = [char]8217   # U+2019 RIGHT SINGLE QUOTATION MARK ('smart' apostrophe)

# --- 1 & nearby context: paragraph 7, "There's a process to making..." ---
$runs7 = @(
    "There${rsq}s a process to making ",
    "decision that ",
    "we${rsq}re probably not thinking",
    " about",
    ". "
)
$null = $d.Paragraphs(7).Range.InsertXML((Build-BodyParagraphXml $runs7))
Write-Host "1) Paragraph 7 -> $($d.Paragraphs(7).Range.Text)"

# --- 2: paragraph 11, "Every time we make a decision..." (drop proofErr) ---
$runs11 = @(
    "Every time we ",
    "make a decision",
    " to do something, we${rsq}re also making ",
    "a decision",
    " not to do ",
    "something else",
    ". ",
    "For any good decision, t",
    "hose alternatives have tradeoffs that need to be considered."
)
$null = $d.Paragraphs(11).Range.InsertXML((Build-BodyParagraphXml $runs11))
Write-Host "2) Paragraph 11 -> $($d.Paragraphs(11).Range.Text)"

# --- 3: paragraph 19, rewrite the "Decisions need to be made objectively..." opening ---
$runs19 = @(
    "Our",
    " emotions and biases",
    " are always going to have an influence on our decisions",
    ". ",
    "For a decision to be",
    " objective",
    ", ",
    "the",
    " ",
    "process",
    "es",
    " that led to the decision need to be (1) explicit and (2) ",
    "consistent. ",
    "AI ensures that these rules are ",
    "met",
    ". "
)
$null = $d.Paragraphs(19).Range.InsertXML((Build-BodyParagraphXml $runs19))
Write-Host "3) Paragraph 19 -> $($d.Paragraphs(19).Range.Text)"

# --- 4 & 5: paragraph 35, Ethics paragraph (drop both proofErr pairs) ---
$runs35 = @(
    "Ethics is why I dove into",
    " Philosophy and Value Theory. Did you know philosophers like Plato, Aristotle, and Aquinas actually provide",
    "d",
    " means to measur",
    "e",
    " ",
    "how ethical ideas are",
    "? If we want to decide whether to invest more in taking care of the sick or feeding the hungry purely on an ethical basis, we actually can. "
)
$null = $d.Paragraphs(35).Range.InsertXML((Build-BodyParagraphXml $runs35))
Write-Host "4/5) Paragraph 35 -> $($d.Paragraphs(35).Range.Text)"

Write-Host "Done. Paragraph count: $($d.Paragraphs.Count)"
